$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 4) with the latest mod-count entry, mirroring
# the existing layout/format used by the previous data row (row 3).

# Column A holds a date written as plain text (e.g. "2025/11/12"), not a
# real date value, so force a text number format first to stop Excel from
# auto-converting the literal into a date serial number, then clear the
# temporary format back off so the cell keeps the sheet's default style.
$dateCell = $ws.Range("A4")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/11/13"
$dateCell.ClearFormats()

$ws.Range("B4").Value = "逃离鸭科夫"
$ws.Range("C4").Value = 1081

# Match the centered alignment used by the rest of the data rows.
$ws.Range("A4:C4").HorizontalAlignment = -4108
$ws.Range("A4:C4").VerticalAlignment = -4108
